$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Summary" to "Validation"
$ws.Name = "Validation"

# Clear the old content (title row + status row are being removed)
$ws.Cells.Clear()

# Rewrite the table: header row, then data rows (no bold title, no status row)
$ws.Range("A1").Value = "Metric"
$ws.Range("B1").Value = "Value"

$ws.Range("A2").Value = "TOC Entries"
$ws.Range("B2").Value = 1360

$ws.Range("A3").Value = "Content Items"
$ws.Range("B3").Value = 25760
